$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.123.59'
$ws.Range("E2").Value = '  +1.59%  '
# Row 3
$ws.Range("D3").Value = '3.141.80'
$ws.Range("E3").Value = '  +3.29%  '
# Row 4
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = '  +0.07%  '
# Row 5
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.46'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +1.56%  '
# Row 6
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.93'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +6.12%  '
# Row 7
$ws.Range("E7").Value = '  -0.08%  '
# Row 8
$ws.Range("D8").Value = '3.134.49'
$ws.Range("E8").Value = '  +3.32%  '
# Row 9
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.501'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +3.15%  '
# Row 10
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.11'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +16.23%  '
# Row 11
$ws.Range("E11").Value = '  +2.36%  '
# Row 12
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  +1.73%  '
# Row 13
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.12'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +2.94%  '
# Row 14
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000224'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +1.36%  '
# Row 15
$ws.Range("D15").Value = '3.647.82'
$ws.Range("E15").Value = '  +3.24%  '
# Row 16
$ws.Range("D16").Value = '65.132.87'
$ws.Range("E16").Value = '  +1.67%  '
# Row 17
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '540.13'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  +12.13%  '
# Row 18
$ws.Range("E18").Value = '  +2.36%  '
# Row 19
$ws.Range("D19").Value = '3.142.32'
$ws.Range("E19").Value = '  +3.16%  '
# Row 20
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.78'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +3.39%  '
# Row 21
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.94'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +3.10%  '
# Row 22
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.709'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +5.40%  '
# Row 23
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.49'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +4.65%  '
# Row 24
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.89'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +3.50%  '
# Row 25
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.07'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +1.69%  '
# Row 26
$ws.Range("E26").Value = '  +0.18%  '
# Row 27
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.90'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +15.70%  '
# Row 28
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.82'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  +3.15%  '
# Row 29
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.15'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +3.52%  '
# Row 30
$ws.Range("E30").Value = '  -0.07%  '
# Row 31
$ws.Range("B31").Value = 'Stacks'
$ws.Range("C31").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.66'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +3.44%  '
# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.31'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +0.76%  '
# Row 33
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.15'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +4.09%  '
# Row 34
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '555.20'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +11.44%  '
# Row 35
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.44'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +1.46%  '
# Row 36
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.11'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  +4.82%  '
# Row 37
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0449'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  +10.21%  '
# Row 38
$ws.Range("E38").Value = '  +0.62%  '
# Row 39
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0822'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +4.84%  '
# Row 40
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.91'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +10.92%  '
# Row 41
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '3.075.51'
$ws.Range("E41").Value = '  +7.91%  '
# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.122'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +3.39%  '
# Row 43
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.30'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +0.62%  '
# Row 44
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.260'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +8.04%  '
# Row 45
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.19'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +7.45%  '
# Row 46
$ws.Range("E46").Value = '  +0.10%  '
# Row 47
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.12'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +2.24%  '
# Row 48
$ws.Range("B48").Value = 'PEPE'
$ws.Range("C48").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D48").Value = '0.0₃0528'
$ws.Range("E48").Value = '  -1.49%  '
# Row 49
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.110'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +2.78%  '
# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '120.00'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +1.76%  '
# Row 51
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.12'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +3.95%  '
